$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Statut" column header in E1, reusing the same header
# formatting (bold, centered, bordered) already applied to the other
# header cells by copying D1's format onto E1.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Statut"

# Match the column width used for the new "Statut" column.
$ws.Columns.Item(5).ColumnWidth = 26.7115

# Update the saved selection/active cell for the sheet.
$ws.Range("E8").Select()
